$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a value as literal text (shared string) without Excel's
# automatic number/date inference changing its type or stamping a new
# number-format style onto the cell. We build the text via a formula in a
# scratch cell (forces a string result), copy it, then paste-special just
# the values into the destination so the destination keeps its existing
# cell formatting/style untouched.
function Set-LiteralText($addr, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range("Z1").Formula = "=""" + $escaped + """"
    $ws.Range("Z1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# Helper: make sure a (possibly brand new) cell ends up with the same
# formatting as the rest of its column by copying formats from a
# known-good neighbour cell in the same column.
function Copy-ColumnFormat($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial(-4122)
}

# Name: "Integrated Project I" -> "Integrated Project"
$ws.Range("B4").Value = "Integrated Project"
$ws.Range("C4").Value = "Integrated Project"

# Créditos-aula: 1 -> 4 (plain numeric-looking text, keep as text)
Set-LiteralText "B5" "4"
Set-LiteralText "C5" "4"

# Carga horária: 135 h -> 180 h
$ws.Range("B7").Value = "180 h"
$ws.Range("C7").Value = "180 h"

# Ativação: 01/01/2012 -> 01/01/2023 (date-looking text, keep as text)
Set-LiteralText "B8" "01/01/2023"
Set-LiteralText "C8" "01/01/2023"

# Objectives: new text added (B11/C11 did not exist before)
$ws.Range("B11").Value = "Introduce students to the principles and methodology of scientific research."
$ws.Range("C11").Value = "Introduce students to the principles and methodology of scientific research."
Copy-ColumnFormat "B10" "B11"
Copy-ColumnFormat "C10" "C11"

# Programa resumido (duplicate date field elsewhere in sheet) also updated
Set-LiteralText "B13" "01/01/2023"
Set-LiteralText "C13" "01/01/2023"

# Short syllabus: new text added (B14/C14 did not exist before)
$ws.Range("B14").Value = "Initiation into a research project under the guidance of a professor."
$ws.Range("C14").Value = "Initiation into a research project under the guidance of a professor."
Copy-ColumnFormat "B13" "B14"
Copy-ColumnFormat "C13" "C14"

# Syllabus: new text added (B16/C16 did not exist before)
$ws.Range("B16").Value = "Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course."
$ws.Range("C16").Value = "Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course."
Copy-ColumnFormat "B15" "B16"
Copy-ColumnFormat "C15" "C16"

# Critério: monografia -> projeto de pesquisa
$ws.Range("B19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."
$ws.Range("C19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de projeto de pesquisa."

# Norma de recuperação: updated text
$ws.Range("B20").Value = "Nota de avaliação do projeto e demais documentos."
$ws.Range("C20").Value = "Nota de avaliação do projeto e demais documentos."

# Bibliografia (recovery-policy text): updated text
$ws.Range("B21").Value = "Devido às características do curso, não será oferecida recuperação."
$ws.Range("C21").Value = "Devido às características do curso, não será oferecida recuperação."

# Clean up the scratch cell used for literal-text writes
$ws.Range("Z1").Clear()
